# Update cryptocurrency price and volume data (refresh snapshot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.916.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.832.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.92%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6872"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.64%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9992"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07652"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3047"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07806"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.832.82"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.070"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6764"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.447"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("E17").Value = "  -1.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.906.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.076.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.420"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1472"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.48%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.770"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.535"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.211"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.105"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.77%  "

$ws.Range("E32").Value = "  -0.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05108"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7535"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.828"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.143"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.674"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01845"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.226.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.693"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.95%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9175"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.74%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9984"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.519"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.62%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5173"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.509"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.975.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.734"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4183"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.81%  "
